$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "RM 232" row (row 26) and the "SC 92" row (originally row 28,
# now row 27 after the first deletion shifts everything up by one).
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()
